$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# "Код реализации" lost its inherited first-line indent from the Normal
# style (it needs a zero hanging indent instead of the style's first-line
# indent), so give its paragraph an explicit (zero) hanging indent.
$rng1 = $d.Content
$rng1.Find.Execute("Код реализации", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$codeParagraph = $rng1.Paragraphs.Item(1)
$codeParagraph.Format.FirstLineIndent = -0.001

# --- Change 2 -------------------------------------------------------------
# Add an empty, centered paragraph (a blank "space") right before
# "Вызов тоста", matching the spacing already used elsewhere before images.
$rng2 = $d.Content
$rng2.Find.Execute("Вызов тоста", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$rng2.InsertParagraphBefore()
$spacerBeforeToast = $rng2.Paragraphs.Item(1)
$spacerBeforeToast.Format.Alignment = 1 # wdAlignParagraphCenter

# --- Change 3 ---------------------------------------------------------
# Add another empty, centered paragraph (a blank "space") right after the
# "Рисунок 4 – вызов Toast" caption, before "Дополнительное задание".
$rng3 = $d.Content
$rng3.Find.Execute("Toast", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$rng3.InsertParagraphAfter()
